# Expense tracker: add a "weekday_text" column (derived from the existing
# numeric weekday column, now renamed "weekday_number") right before the
# "store" column, so the sheet reads: date, expense_category, expense_type,
# value, month, year, weekday_number, weekday_text, store, city,
# english_translation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (pushes store/city/english_translation right by one).
$ws.Columns("H").Insert()

# Rename the existing WEEKDAY() header and label the new column.
$ws.Range("G1").Value = "weekday_number"
$ws.Range("H1").Value = "weekday_text"

# First data row gets its own (non-shared) formula...
$ws.Range("H2").Formula = '=CHOOSE(WEEKDAY(A2, 2), "Monday", "Tuesday","Wednesday", "Thursday", "Friday", "Saturday","Sunday")'

# ...the rest of the rows are filled as a shared formula block.
$ws.Range("H3:H6").Formula = '=CHOOSE(WEEKDAY(A3, 2), "Monday", "Tuesday","Wednesday", "Thursday", "Friday", "Saturday","Sunday")'

# Keep the hidden filter-database defined name in sync with the new column.
$fd = $wb.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "=total!`$A`$1:`$K`$6"

# Match the recorded cursor position left behind by the edit.
[void]$ws.Range("H13").Select()
